$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post corresponding to row 401 ("「彼の尖った牙を心配するな」...") was removed.
# Deleting the entire row shifts all subsequent rows up by one, matching the diff.
$ws.Rows.Item(401).EntireRow.Delete()
